$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 2.0938545198362832
$ws.Cells.Item(2, 1).Value = 2.0953089440720643
$ws.Cells.Item(3, 1).Value = 2.0952132171374087
$ws.Cells.Item(4, 1).Value = 2.095166371016555
$ws.Cells.Item(5, 1).Value = 2.1000727170443314
$ws.Cells.Item(6, 1).Value = 2.1045698923043954
$ws.Cells.Item(7, 1).Value = 2.1146161073894323
$ws.Cells.Item(8, 1).Value = 2.1192481160072463
$ws.Cells.Item(9, 1).Value = 2.1203623492637274
$ws.Cells.Item(10, 1).Value = 2.133027091882588
$ws.Cells.Item(11, 1).Value = 2.142250822166977
$ws.Cells.Item(12, 1).Value = 2.1394085053083467
$ws.Cells.Item(13, 1).Value = 2.1402803176651046
$ws.Cells.Item(14, 1).Value = 2.1426829645685803
$ws.Cells.Item(15, 1).Value = 2.1404454009210823
$ws.Cells.Item(16, 1).Value = 2.1398118071610313
$ws.Cells.Item(17, 1).Value = 2.1395820986637055
$ws.Cells.Item(18, 1).Value = 2.140898636735457
$ws.Cells.Item(19, 1).Value = 2.145026321329852
$ws.Cells.Item(20, 1).Value = 2.1457064397498415
$ws.Cells.Item(21, 1).Value = 2.149524483892052
$ws.Cells.Item(22, 1).Value = 2.15416376313208
$ws.Cells.Item(23, 1).Value = 2.153145859182122
$ws.Cells.Item(24, 1).Value = 2.1218769733475242
$ws.Cells.Item(25, 1).Value = 2.0987079486242557
$ws.Cells.Item(26, 1).Value = 2.0687593823665296
$ws.Cells.Item(27, 1).Value = 2.047623063176421
$ws.Cells.Item(28, 1).Value = 2.037779597025147
$ws.Cells.Item(29, 1).Value = 2.0241554013966736
$ws.Cells.Item(30, 1).Value = 2.019593051730899
$ws.Cells.Item(31, 1).Value = 2.0108078116646677
$ws.Cells.Item(32, 1).Value = 2.0020249891198767
$ws.Cells.Item(33, 1).Value = 2.002322869550157
$ws.Cells.Item(34, 1).Value = 1.9944574228583967
$ws.Cells.Item(35, 1).Value = 1.9915366956113907
$ws.Cells.Item(36, 1).Value = 1.990341603175589
$ws.Cells.Item(37, 1).Value = 1.9871487509794927
$ws.Cells.Item(38, 1).Value = 1.9832747415817984
$ws.Cells.Item(39, 1).Value = 1.979431747232117
$ws.Cells.Item(40, 1).Value = 1.9719463783988673
$ws.Cells.Item(41, 1).Value = 1.9668676793937314
$ws.Cells.Item(42, 1).Value = 1.962372473330539
$ws.Cells.Item(43, 1).Value = 1.9595395452632034
$ws.Cells.Item(44, 1).Value = 1.949764103949232
$ws.Cells.Item(45, 1).Value = 1.9403791414250815
$ws.Cells.Item(46, 1).Value = 1.9236781223863049
$ws.Cells.Item(47, 1).Value = 1.91384554902122
$ws.Cells.Item(48, 1).Value = 1.903935124246722
$ws.Cells.Item(49, 1).Value = 1.9043305130964843
$ws.Cells.Item(50, 1).Value = 1.9009231581589865
$ws.Cells.Item(51, 1).Value = 1.8846319375270713
$ws.Cells.Item(52, 1).Value = 1.872814416465315
$ws.Cells.Item(53, 1).Value = 1.8713374194754464
$ws.Cells.Item(54, 1).Value = 1.8725342800430325
$ws.Cells.Item(55, 1).Value = 1.870593866901955
$ws.Cells.Item(56, 1).Value = 1.8673605676835132
$ws.Cells.Item(57, 1).Value = 1.867307722341149
$ws.Cells.Item(58, 1).Value = 1.8678483088994273
$ws.Cells.Item(59, 1).Value = 1.8669147125337506
$ws.Cells.Item(60, 1).Value = 1.862319754909743
$ws.Cells.Item(61, 1).Value = 1.8561606897962803
$ws.Cells.Item(62, 1).Value = 1.853801956628585
$ws.Cells.Item(63, 1).Value = 1.8541894139770312
$ws.Cells.Item(64, 1).Value = 1.8495927098437392
$ws.Cells.Item(65, 1).Value = 1.8471097327458272
$ws.Cells.Item(66, 1).Value = 1.8376911241497689
$ws.Cells.Item(67, 1).Value = 1.8364830029601324
$ws.Cells.Item(68, 1).Value = 1.834157325448686
$ws.Cells.Item(69, 1).Value = 1.833082987340787
$ws.Cells.Item(70, 1).Value = 1.8322535634633295
$ws.Cells.Item(71, 1).Value = 1.8331000461612605
$ws.Cells.Item(72, 1).Value = 1.832264546293851
$ws.Cells.Item(73, 1).Value = 1.8349575561311824
$ws.Cells.Item(74, 1).Value = 1.8340081866653775
$ws.Cells.Item(75, 1).Value = 1.8340003924954633
$ws.Cells.Item(76, 1).Value = 1.8348410152354835
$ws.Cells.Item(77, 1).Value = 1.8349479168532656
$ws.Cells.Item(78, 1).Value = 1.8323964738699114
$ws.Cells.Item(79, 1).Value = 1.8253721080588474
$ws.Cells.Item(80, 1).Value = 1.8254717916675276
$ws.Cells.Item(81, 1).Value = 1.8206538427081256
$ws.Cells.Item(82, 1).Value = 1.8194103109856186
$ws.Cells.Item(83, 1).Value = 1.8128857400766676
$ws.Cells.Item(84, 1).Value = 1.8018411271547308
$ws.Cells.Item(85, 1).Value = 1.7889963834082119
$ws.Cells.Item(86, 1).Value = 1.777525097779157
$ws.Cells.Item(87, 1).Value = 1.7679011812062997
$ws.Cells.Item(88, 1).Value = 1.764652866740293
$ws.Cells.Item(89, 1).Value = 1.758154451405913
$ws.Cells.Item(90, 1).Value = 1.755158034566929
$ws.Cells.Item(91, 1).Value = 1.753154231699714
$ws.Cells.Item(92, 1).Value = 1.752828493524194
$ws.Cells.Item(93, 1).Value = 1.7529571475689387
$ws.Cells.Item(94, 1).Value = 1.751292787048632
$ws.Cells.Item(95, 1).Value = 1.748650406280066
$ws.Cells.Item(96, 1).Value = 1.7466178012258649
$ws.Cells.Item(97, 1).Value = 1.7437544978205834
$ws.Cells.Item(98, 1).Value = 1.74462884783254
$ws.Cells.Item(99, 1).Value = 1.7473730952324922
$ws.Cells.Item(100, 1).Value = 1.7503259712375363
$ws.Cells.Item(101, 1).Value = 1.7557850003540887
$ws.Cells.Item(102, 1).Value = 1.7583625952496396
$ws.Cells.Item(103, 1).Value = 1.7597317523466853
$ws.Cells.Item(104, 1).Value = 1.758139193764323
$ws.Cells.Item(105, 1).Value = 1.7536077664999348
$ws.Cells.Item(106, 1).Value = 1.7481001041456627
$ws.Cells.Item(107, 1).Value = 1.740433437793139
$ws.Cells.Item(108, 1).Value = 1.736446535782509
$ws.Cells.Item(109, 1).Value = 1.7326676209824785
$ws.Cells.Item(110, 1).Value = 1.7312759578318888
$ws.Cells.Item(111, 1).Value = 1.7319849094880864
$ws.Cells.Item(112, 1).Value = 1.7316930070653138
$ws.Cells.Item(113, 1).Value = 1.7293276326678977
$ws.Cells.Item(114, 1).Value = 1.7268542136818685
$ws.Cells.Item(115, 1).Value = 1.7265385736611163
$ws.Cells.Item(116, 1).Value = 1.726653158503806
$ws.Cells.Item(117, 1).Value = 1.72834037446476
$ws.Cells.Item(118, 1).Value = 1.7258180756839994
$ws.Cells.Item(119, 1).Value = 1.7249687554495585
$ws.Cells.Item(120, 1).Value = 1.720665180414093
$ws.Cells.Item(121, 1).Value = 1.7106063700738527
$ws.Cells.Item(122, 1).Value = 1.697005250534131
$ws.Cells.Item(123, 1).Value = 1.6867839868400618
$ws.Cells.Item(124, 1).Value = 1.6746182092525537
$ws.Cells.Item(125, 1).Value = 1.6816215734263829
$ws.Cells.Item(126, 1).Value = 1.6808658844190094
$ws.Cells.Item(127, 1).Value = 1.6639059801623055
$ws.Cells.Item(128, 1).Value = 1.6548497992016187
$ws.Cells.Item(129, 1).Value = 1.639446113283297
$ws.Cells.Item(130, 1).Value = 1.629023431149009
$ws.Cells.Item(131, 1).Value = 1.6258151435416934
$ws.Cells.Item(132, 1).Value = 1.6279573970061794
$ws.Cells.Item(133, 1).Value = 1.6297294242422118
$ws.Cells.Item(134, 1).Value = 1.6320582791733584
$ws.Cells.Item(135, 1).Value = 1.630032070997752
$ws.Cells.Item(136, 1).Value = 1.6294883446988475
$ws.Cells.Item(137, 1).Value = 1.6253259396747275
$ws.Cells.Item(138, 1).Value = 1.6243167304515018
$ws.Cells.Item(139, 1).Value = 1.6304933839772833
$ws.Cells.Item(140, 1).Value = 1.641041111759645
$ws.Cells.Item(141, 1).Value = 1.6462527644210847
$ws.Cells.Item(142, 1).Value = 1.6470733680034346
$ws.Cells.Item(143, 1).Value = 1.6437175503238897
$ws.Cells.Item(144, 1).Value = 1.6360315853797303
$ws.Cells.Item(145, 1).Value = 1.6288506002655945
$ws.Cells.Item(146, 1).Value = 1.6238766101665965
$ws.Cells.Item(147, 1).Value = 1.6248292140192788
$ws.Cells.Item(148, 1).Value = 1.6328124567124593
$ws.Cells.Item(149, 1).Value = 1.6500501631942144
$ws.Cells.Item(150, 1).Value = 1.6646029473024901
$ws.Cells.Item(151, 1).Value = 1.6659461871719203
$ws.Cells.Item(152, 1).Value = 1.6719379286375256
$ws.Cells.Item(153, 1).Value = 1.6706067181340538
$ws.Cells.Item(154, 1).Value = 1.6647187346234313
$ws.Cells.Item(155, 1).Value = 1.6656434144829033
$ws.Cells.Item(156, 1).Value = 1.6476484734768477
$ws.Cells.Item(157, 1).Value = 1.6403767696208993
$ws.Cells.Item(158, 1).Value = 1.6311517200037753
$ws.Cells.Item(159, 1).Value = 1.6199140843629576
$ws.Cells.Item(160, 1).Value = 1.619510663368561
$ws.Cells.Item(161, 1).Value = 1.6158103489939801
$ws.Cells.Item(162, 1).Value = 1.6158406209360943
$ws.Cells.Item(163, 1).Value = 1.6242073694631975
$ws.Cells.Item(164, 1).Value = 1.6327374467667801
$ws.Cells.Item(165, 1).Value = 1.6427488374433845
$ws.Cells.Item(166, 1).Value = 1.644165346847282
$ws.Cells.Item(167, 1).Value = 1.6459305796106567
$ws.Cells.Item(168, 1).Value = 1.6395130679189813
$ws.Cells.Item(169, 1).Value = 1.6334396028262517
$ws.Cells.Item(170, 1).Value = 1.6337602658366737
$ws.Cells.Item(171, 1).Value = 1.6301523751364182
$ws.Cells.Item(172, 1).Value = 1.631723668212752
$ws.Cells.Item(173, 1).Value = 1.6326613188192711
$ws.Cells.Item(174, 1).Value = 1.6325403036245008
$ws.Cells.Item(175, 1).Value = 1.6348578179593576
$ws.Cells.Item(176, 1).Value = 1.6386886733477004
$ws.Cells.Item(177, 1).Value = 1.6405812992472537
$ws.Cells.Item(178, 1).Value = 1.6425319617956915
$ws.Cells.Item(179, 1).Value = 1.6446391166690535
$ws.Cells.Item(180, 1).Value = 1.6481277964637497
$ws.Cells.Item(181, 1).Value = 1.6462293273122697
$ws.Cells.Item(182, 1).Value = 1.6495084494252483
$ws.Cells.Item(183, 1).Value = 1.6489351895714157
$ws.Cells.Item(184, 1).Value = 1.6521139809905458
$ws.Cells.Item(185, 1).Value = 1.6567971045445322
$ws.Cells.Item(186, 1).Value = 1.657418993368565
$ws.Cells.Item(187, 1).Value = 1.6624798186319119
$ws.Cells.Item(188, 1).Value = 1.664906996512153
$ws.Cells.Item(189, 1).Value = 1.664357189814254
$ws.Cells.Item(190, 1).Value = 1.6651152246033831
$ws.Cells.Item(191, 1).Value = 1.664352156745647
$ws.Cells.Item(192, 1).Value = 1.6647311928467585
$ws.Cells.Item(193, 1).Value = 1.6642905775694419
$ws.Cells.Item(194, 1).Value = 1.6607533399383954
$ws.Cells.Item(195, 1).Value = 1.650313756460619
$ws.Cells.Item(196, 1).Value = 1.6483009059722376
$ws.Cells.Item(197, 1).Value = 1.6487902964613834
$ws.Cells.Item(198, 1).Value = 1.647796605261629
$ws.Cells.Item(199, 1).Value = 1.6466498052821998
$ws.Cells.Item(200, 1).Value = 1.646748287573634

Write-Host "Updated 200 cells in column A"
